$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to re-pulled/recalculated data
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -2
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = -5
$ws.Range("F17").Value = -3
$ws.Range("F21").Value = 1
$ws.Range("F24").Value = -5
$ws.Range("F30").Value = -3
$ws.Range("F34").Value = 0
